$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Coin/Link/Price/Volume columns (B:E) hold plain text values in this
# sheet -- including price strings that use "." as a thousands separator
# (e.g. "69.145.66") which Excel would otherwise coerce to a number when
# assigned via .Value. Temporarily force a text number format while writing
# the refreshed figures, then restore the original cell style so formatting
# is left exactly as it was.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$updates = @(
    ,@('D2', '69.145.66')
    ,@('E2', '  -2.46%  ')
    ,@('D3', '3.690.09')
    ,@('E3', '  -3.29%  ')
    ,@('E4', '  -0.12%  ')
    ,@('D5', '680.29')
    ,@('E5', '  -3.84%  ')
    ,@('D6', '162.51')
    ,@('E6', '  -4.40%  ')
    ,@('D7', '3.690.73')
    ,@('E7', '  -3.25%  ')
    ,@('D8', '0.999')
    ,@('E8', '  -0.11%  ')
    ,@('D9', '0.492')
    ,@('E9', '  -5.70%  ')
    ,@('D10', '0.148')
    ,@('E10', '  -7.58%  ')
    ,@('D11', '7.19')
    ,@('E11', '  -2.27%  ')
    ,@('D12', '0.449')
    ,@('E12', '  -1.48%  ')
    ,@('D13', '0.0000235')
    ,@('E13', '  -7.12%  ')
    ,@('D14', '33.42')
    ,@('E14', '  -7.92%  ')
    ,@('D15', '4.304.75')
    ,@('E15', '  -3.55%  ')
    ,@('D16', '3.677.10')
    ,@('E16', '  -3.70%  ')
    ,@('D17', '69.213.95')
    ,@('E17', '  -2.39%  ')
    ,@('D19', '16.36')
    ,@('E19', '  -5.30%  ')
    ,@('D20', '6.58')
    ,@('E20', '  -8.19%  ')
    ,@('D21', '482.60')
    ,@('E21', '  -2.15%  ')
    ,@('D22', '9.82')
    ,@('E22', '  -7.13%  ')
    ,@('D23', '0.663')
    ,@('E23', '  -8.95%  ')
    ,@('D24', '79.26')
    ,@('E24', '  -7.39%  ')
    ,@('D25', '3.826.33')
    ,@('E25', '  -3.68%  ')
    ,@('D26', '0.0000129')
    ,@('E26', '  -10.48%  ')
    ,@('D27', '11.55')
    ,@('E27', '  -4.19%  ')
    ,@('E28', '  +0.04%  ')
    ,@('D29', '9.49')
    ,@('E29', '  -9.93%  ')
    ,@('D30', '1.81')
    ,@('E30', '  -12.03%  ')
    ,@('D31', '2.73')
    ,@('E31', '  -11.54%  ')
    ,@('D32', '2.10')
    ,@('E32', '  -5.29%  ')
    ,@('D33', '6.73')
    ,@('E33', '  -8.64%  ')
    ,@('D34', '1.00')
    ,@('E34', '  -0.04%  ')
    ,@('D35', '26.73')
    ,@('E35', '  -8.37%  ')
    ,@('D36', '0.164')
    ,@('E36', '  -5.61%  ')
    ,@('D37', '3.649.84')
    ,@('E37', '  -3.63%  ')
    ,@('D38', '8.51')
    ,@('E38', '  -6.30%  ')
    ,@('D39', '6.05')
    ,@('E39', '  +2.05%  ')
    ,@('D40', '0.0932')
    ,@('E40', '  -8.18%  ')
    ,@('B41', 'Stacks')
    ,@('C41', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx')
    ,@('D41', '2.20')
    ,@('E41', '  -4.63%  ')
    ,@('B42', 'USDe')
    ,@('C42', 'https://coinranking.com/coin/exbfr2U-0+usde-usde')
    ,@('D42', '1.00')
    ,@('E42', '  -0.01%  ')
    ,@('D43', '0.999')
    ,@('E43', '  -0.20%  ')
    ,@('D44', '0.956')
    ,@('E44', '  -8.33%  ')
    ,@('D45', '160.09')
    ,@('E45', '  -2.31%  ')
    ,@('D46', '48.32')
    ,@('E46', '  -0.89%  ')
    ,@('D47', '2.84')
    ,@('E47', '  -13.36%  ')
    ,@('B48', 'Bittensor')
    ,@('C48', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao')
    ,@('D48', '396.87')
    ,@('E48', '  -6.76%  ')
    ,@('B49', 'ONDO')
    ,@('C49', 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo')
    ,@('D49', '1.30')
    ,@('E49', '  -3.77%  ')
    ,@('B50', 'InjectiveProtocol')
    ,@('C50', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj')
    ,@('D50', '28.39')
    ,@('E50', '  +2.17%  ')
    ,@('B51', 'FLOKI')
    ,@('C51', 'https://coinranking.com/coin/fmHk13Rqw+floki-floki')
    ,@('D51', '0.000275')
    ,@('E51', '  -11.02%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$dataRange.Style = $origStyle
